# feat: support struct with fields spanning multiple columns
#
# The "Activity" sheet gains three new trailing columns (O, P, Q) that hold a
# struct-like value (task type / task param / task target) spanning multiple
# columns, for every data row that previously had content in columns A:N.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Item")
$ws2 = $wb.Worksheets.Item("Activity")

# --- New headers in row 1 (O1, P1, Q1) -----------------------------------
# Written in this order so the new shared-string table entries line up the
# same way the original author's edit produced them.
$ws2.Range("O1").Value = "任务类型"
$ws2.Range("Q1").Value = "任务目标"
$ws2.Range("P1").Value = "任务参数"

# --- New data cells for the populated rows (2,3,4,6,7,8) ------------------
# Row 5 is left untouched -- it never had data beyond column E, so it keeps
# its original shape (spans 1:14).
$ws2.Range("O2").Value = 1
$ws2.Range("P2").Value = 5
$ws2.Range("Q2").Value = 1

$ws2.Range("O3").Value = 1
$ws2.Range("P3").Value = 5
$ws2.Range("Q3").Value = 2

$ws2.Range("O4").Value = 1
$ws2.Range("P4").Value = 5
$ws2.Range("Q4").Value = 3

$ws2.Range("O6").Value = 2
$ws2.Range("P6").Value = 7
$ws2.Range("Q6").Value = 1

$ws2.Range("O7").Value = 2
$ws2.Range("P7").Value = 7
$ws2.Range("Q7").Value = 2

$ws2.Range("O8").Value = 2
$ws2.Range("P8").Value = 7
$ws2.Range("Q8").Value = 3

# --- Selection / active-tab bookkeeping -----------------------------------
# The author ended up with the Activity sheet active (having just typed the
# new column data), with the cursor resting on P12, while the Item sheet's
# previous selection moved off N10 onto G11.
[void]$ws1.Range("G11").Select()
[void]$ws2.Activate()
[void]$ws2.Range("P12").Select()
